$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ShipmentTracking numbers (column P) for rows 2-25, and the updated
# ActualRate (column Q) for row 22 -- "cheetah processing" re-run assigned
# fresh FedEx tracking numbers / rates to the existing shipment rows.
$trackingUpdates = @{
    2  = '320018179991'
    3  = '320018180002'
    4  = '320018180035'
    5  = '320018180057'
    6  = '320018180090'
    7  = '320018180127'
    8  = '320018180150'
    9  = '320018180171'
    10 = '320018180208'
    11 = '320018180220'
    12 = '320018180263'
    13 = '320018180285'
    14 = '320018180311'
    15 = '320018180333'
    16 = '320018180366'
    17 = '320018180388'
    18 = '320018180425'
    19 = '320018180447'
    20 = '320018180480'
    21 = '320018180506'
    22 = '320018180539'
    23 = '320018180540'
    24 = '320018180550'
    25 = '320018180561'
}

foreach ($row in $trackingUpdates.Keys) {
    # Leading apostrophe forces these digit-only strings to stay text
    # (shared-string) cells instead of being parsed as numbers.
    $ws.Range("P$row").Value = "'" + $trackingUpdates[$row]
}

# Row 22's ActualRate also changed as part of this re-run.
$ws.Range("Q22").Value = "'`$202.67"

# Re-normalize style back to the default (no explicit style id), matching
# the original unstyled P/Q cells -- the text-coercion above stamps a
# quote-prefixed "Text" style that we don't want to keep.
$ws.Range("P2:P25").Style = "Normal"
$ws.Range("Q22").Style = "Normal"
